$d = $word.ActiveDocument

# --- Change 1 ---
# The paragraph "Ocorrem muitas falhas de software..." was originally split
# into two runs ("...atende" / "r mais...") with a _GoBack bookmark between
# them. The edit merges the text back into a single run and removes the
# bookmark. A plain Find/Execute replacement (which rewrites the whole
# matched span as one run) reproduces this exactly.
$rng1 = $d.Content
$rng1.Find.Execute(
    "atender mais",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "atender mais", 2
) | Out-Null

# --- Change 2 ---
# "conforma a demanda" becomes "conforme a demanda", but the resulting text
# is split into three runs: "...conform", "e", then the _GoBack bookmark,
# then " a demanda...". This mirrors someone selecting the trailing "a" and
# retyping "e" (leaving the caret - and therefore _GoBack - right after it).
$rng2 = $d.Content
$rng2.Find.Execute(
    "conforma",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "", 0
) | Out-Null

$wordEnd = $rng2.End

# Replace the trailing "a" of "conforma" with "e".
$charRng = $d.Range($wordEnd - 1, $wordEnd)
$charRng.Text = "e"

# Force a run boundary between "conform" and "e" by briefly bookmarking the
# seam, then discard that helper bookmark (the split it produced persists).
$seam = $d.Range($wordEnd - 1, $wordEnd - 1)
$d.Bookmarks.Add("TempSplit", $seam) | Out-Null
$d.Bookmarks("TempSplit").Delete()

# Re-create the real _GoBack bookmark right after the newly typed "e".
$goBack = $d.Range($wordEnd, $wordEnd)
$d.Bookmarks.Add("_GoBack", $goBack) | Out-Null
